$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column B for nameId|String (shifts star, accumulatedAtk, prefabAddress right)
$ws.Columns("B").Insert()

# Insert new column F for spriteName|String (after prefabAddress, which is now column E)
$ws.Columns("F").Insert()

# Header row - enter new header text in the same order Excel originally recorded
# them in sharedStrings.xml (spriteName|String, then nameId|String)
$ws.Range("F1").Value = "spriteName|String"
$ws.Range("B1").Value = "nameId|String"
$ws.Range("A1").Value = "petId|String"
$ws.Range("C1").Value = "star|Int"
$ws.Range("D1").Value = "accumulatedAtk|Int"
$ws.Range("E1").Value = "prefabAddress|String"
$ws.Range("G1").Value = "orderIndex|Int"
$ws.Range("H1").Value = "meetWeight|Float"

# Data rows - PetPortrait_0001 first, then the PetName_000x series
$ws.Range("F2").Value = "PetPortrait_0001"
$ws.Range("B2").Value = "PetName_0001"
$ws.Range("B3").Value = "PetName_0002"
$ws.Range("B4").Value = "PetName_0003"
$ws.Range("B5").Value = "PetName_0004"
$ws.Range("B6").Value = "PetName_0005"
$ws.Range("B7").Value = "PetName_0006"

$ws.Range("F3:F7").Value = "PetPortrait_0001"

# Column widths (target XML widths are 14.25 and 16 "characters"; the
# character->pixel->XML round trip this host applies only lands on
# multiples of 1/7, so we dial in the closest input that reproduces the
# desired stored width as exactly as possible)
$ws.Columns("B").ColumnWidth = 13.5
$ws.Columns("F").ColumnWidth = 15.28
